# Redo the calcs in the spreadsheet: duplicate Sheet1 into a new "Sheet2"
# with reworked LED layout calculations, keep the original as "Pre-resize".

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Deselect the old explicit active-cell on the original sheet and select
# the first two columns instead (matches the recorded view state).
$ws1.Range("A1:B1048576").Select()

# Duplicate the sheet (this carries over column widths/row heights/styles)
# and use the copy as the new "Sheet2" that holds the reworked numbers.
$ws1.Copy($null, $ws1)
$ws2 = $wb.Worksheets.Item(2)

# Rename sheets: original -> "Pre-resize", duplicate -> "Sheet2"
$ws1.Name = "Pre-resize"
$ws2.Name = "Sheet2"

# Wipe all of the copied content/formatting so we can rebuild it.
$ws2.Cells.Clear()

# --- Rebuild Sheet2 content -------------------------------------------

# X / Y
$ws2.Range("A1").Value = "X"
$ws2.Range("B1").Value = 530

$ws2.Range("A2").Value = "Y"
$ws2.Range("B2").Value = 401

# Aspect ratio
$ws2.Range("A3").Value = "Aspect ratio"
$ws2.Range("B3").Formula = "=B1/B2"

# Total pixels
$ws2.Range("A4").Value = "Total pixels"
$ws2.Range("B4").Formula = "=B1*B2"

# Pixels in 0 / Pixels in 1
$ws2.Range("A6").Value = "Pixels in 0"
$ws2.Range("B6").Value = 80385

$ws2.Range("A7").Value = "Pixels in 1"
$ws2.Range("B7").Value = 36492

# Pixel ratio
$ws2.Range("A8").Value = "Pixel ratio"
$ws2.Range("B8").Formula = "=(B6+B7)/B4"

# LED spacing (new row/label)
$ws2.Range("A10").Value = "LED spacing"
$ws2.Range("B10").Value = 33.3333

# LED area (mm2)
$ws2.Range("A11").Value = "LED area (mm2)"
$ws2.Range("B11").Formula = "=B10*B10"

# Header row for the results table
$ws2.Range("B13").Value = "Height"
$ws2.Range("C13").Value = "Width"
$ws2.Range("D13").Value = "Y"
$ws2.Range("E13").Value = "X"
$ws2.Range("F13").Value = "LED count"

# Results row
$ws2.Range("A14").Value = "Option 1"
$ws2.Range("B14").Value = 985
$ws2.Range("C14").Formula = "=B14*B3"
$ws2.Range("D14").Formula = "=ROUNDDOWN(B14/B10,0)"
$ws2.Range("E14").Formula = "=ROUNDDOWN(C14/B10,0)"
$ws2.Range("F14").Formula = "=D14*E14*B8"

# --- Styling -------------------------------------------------------------
# Column A keeps the bold font inherited from the original sheet; make it
# left aligned (matches the "LED spacing"/"Height"/etc label column look).
$ws2.Range("A1:A14").HorizontalAlignment = -4131

# Data columns (B:F) get left alignment too.
$ws2.Range("B1:F14").HorizontalAlignment = -4131

# Select the last computed cell, matching the recorded view state.
$ws2.Range("F14").Select()
